# Applies the CPSprint1.docx Sprint-1 edit:
#   1. Appends " (no tiene ediciones registrada" + "s)" (split across four
#      runs, exactly as Word would leave them after an interrupted typing
#      session) to the "Consultar Ediciones de un Torneo" test-case title
#      that sits in the row labelled "4" (the second occurrence of that
#      phrase in the document).
#   2. Moves the "_GoBack" bookmark so it sits between "...registrada" and
#      "s)" - i.e. right where the user's cursor was left - and removes it
#      from its old home in the empty paragraph just after the table.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-PkgXml([string]$bodyXml) {
    return '<?xml version="1.0" encoding="utf-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document ' + $wNs + '><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData>' +
        '</pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# Step 1: locate the SECOND "Consultar Ediciones de un Torneo" (the one
# in the "4" test-case row) and rebuild its paragraph with the new runs
# plus the relocated _GoBack bookmark.
# ---------------------------------------------------------------------
$needle = "Consultar Ediciones de un Torneo"

$find = $d.Content
$find.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$find.Collapse(0)
$find.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$targetPara = $find.Paragraphs(1).Range

$runRpr = '<w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:b/><w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/></w:rPr>'

$newParaXml =
    '<w:p w:rsidR="008B60B0" w:rsidRPr="008B60B0" w:rsidRDefault="005A70B4" w:rsidP="004A21FC">' +
        '<w:pPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/></w:rPr></w:pPr>' +
        '<w:r w:rsidRPr="00987262">' + $runRpr + '<w:t>Consultar Ediciones de un Torneo</w:t></w:r>' +
        '<w:r>' + $runRpr + '<w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:r>' + $runRpr + '<w:t>(no tiene</w:t></w:r>' +
        '<w:r>' + $runRpr + '<w:t xml:space="preserve"> ediciones registrada</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
        '<w:r>' + $runRpr + '<w:t>s)</w:t></w:r>' +
    '</w:p>'

$targetPara.InsertXML((New-PkgXml $newParaXml)) | Out-Null

# ---------------------------------------------------------------------
# Step 2: the old _GoBack bookmark lived by itself in the empty paragraph
# immediately following the table; clear it out so that paragraph goes
# back to being a plain empty paragraph.
# ---------------------------------------------------------------------
$tbl = $d.Tables(1)
$afterTable = $tbl.Range.End
$oldBookmarkPara = $d.Range($afterTable, $afterTable).Paragraphs(1).Range

$oldBookmarkPara.InsertXML((New-PkgXml '<w:p/>')) | Out-Null

Write-Output "Edit applied."
